$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new error code rows, mirroring the style of existing rows (2 and 3)
$ws.Range("A4").Value = 412202
$ws.Range("B4").Value = "SubscriptionLimitExceeded - The user already owns the maximum allowed number of subscriptions."

$ws.Range("A5").Value = 412203
$ws.Range("B5").Value = "TopicLimitExceeded - The user already owns the maximum allowed number of topics."

# Copy the style (left-aligned number format) from A3 to the new A4/A5 cells
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4:A5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Update selection to match the new active cell B5
$ws.Range("B5").Select() | Out-Null
